# Generate Report for Handback
# Refresh the localization-status report: the ca6eecfb-577d-4341-b61f-f8ac843d3b80
# handback for zh-cn / de-de finished successfully (it previously showed as
# "Ready for handoff" / stale-handback error). Update the Overview rollup and
# each language sheet to reflect the new "Handed back: in sync with en-US"
# status plus the fresh handback timestamps, and clear the now-stale error
# detail column.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the ca6eecfb-... file ---
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is the ca6eecfb-... file ---
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("K3").Value = "2016-09-02 00:59:37"
$ws2.Range("P3").Value = ""

# --- de-de sheet: row 3 is the ca6eecfb-... file ---
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("K3").Value = "2016-09-02 00:59:44"
$ws3.Range("P3").Value = ""

# The "Error Detail" column no longer holds a long URL-laden message, so its
# display width collapses back down from the 40-char placeholder.
$ws2.Columns.Item(16).ColumnWidth = 12.83
$ws3.Columns.Item(16).ColumnWidth = 12.83
